# Apply odds/score updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.3
$ws.Range("X2").Value = 1.22

# Row 3
$ws.Range("N3").Value = 8

# Row 4
$ws.Range("G4").Value = 2.7
$ws.Range("I4").Value = 2.55
$ws.Range("J4").Value = 3.5
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("S4").Value = 2.25
$ws.Range("T4").Value = 1.62
$ws.Range("AA4").Value = 1.87
$ws.Range("AB4").Value = 1.77
$ws.Range("AC4").Value = 7.5
$ws.Range("AF4").Value = 29
$ws.Range("AK4").Value = 17
$ws.Range("AM4").Value = 7.5
$ws.Range("AS4").Value = 351

# Row 5
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 3.6
$ws.Range("J5").Value = 2.75
$ws.Range("L5").Value = 4.33
$ws.Range("AN5").Value = 17
$ws.Range("AQ5").Value = 29
